# Express import template: rework header labels + sample row, switch several
# text-like columns (Date, Invoice) to literal text, and add a blank staging
# row beneath the sample so the importer has a ready second data line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels -------------------------------------------------
# All header cells are formatted as Text before the label is typed in so
# re-running this script (or retyping) never lets Excel reinterpret a label.
$ws.Range("A1:G1").NumberFormat = "@"
$ws.Range("A1").Value = "Dept"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Supplier"
$ws.Range("D1").Value = "Invoice"
$ws.Range("E1").Value = "Code"
$ws.Range("F1").Value = "Qty"
$ws.Range("G1").Value = "UnitCost"

# --- Row 2: sample data line ----------------------------------------------
# A2 (Dept) stays as-is ("MBK"); format it to match the rest of the text run.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "MBK"

# B2 (Date) becomes a literal text value instead of a real date serial -
# format as Text first so the "01/10/68" string isn't parsed as a date.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "01/10/68"

# C2 (Supplier) unchanged content, now plain text formatted like its peers.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "026959000"

# D2 (Invoice) switches from a numeric BillNumber to a literal text invoice
# number - format as Text first so the digit string isn't parsed as a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "6500116099"

# E2 (Code) stays text "001".
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "001"

# F2 (Qty) / G2 (UnitCost) remain real numbers - enter the value first, then
# apply the (cosmetic) Text number format so the stored value stays numeric.
$ws.Range("F2").Value = 1
$ws.Range("F2").NumberFormat = "@"

$ws.Range("G2").Value = 9910.57
$ws.Range("G2").NumberFormat = "@"

# Give E2/F2/G2 the "Comma" based text style (matches D2's former Comma
# parentage) so they share one style group distinct from the plain A2:D2 run.
$ws.Range("E2:G2").Style = "Comma"
$ws.Range("E2:G2").NumberFormat = "@"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9910.57

# --- Row 3: blank staging row ----------------------------------------------
# Pre-format E3:G3 as Text/Comma (no values yet) so the row exists, ready for
# the next import line, matching the style used by E2:G2.
$ws.Range("E3:G3").Style = "Comma"
$ws.Range("E3:G3").NumberFormat = "@"

# Put the selection where the user would resume typing.
$ws.Range("D3").Select

Write-Output "express_import_template updated"
